# Add two new columns I (I0) and J (IF) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy the style used by the existing header cells (e.g. H1)
# onto the new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-15
$data = @{
    2  = @(1, 4)
    3  = @(1, 4)
    4  = @(4, 9)
    5  = @(1, 4)
    6  = @(1, 5)
    7  = @(1, 3)
    8  = @(1, 8)
    9  = @(1, 6)
    10 = @(1, 4)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 3)
    15 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
